# Coordenadas.xlsx update:
# Luffy now has a symmetric key + encrypted coordinates token (previously
# that pair lived on Sanji's and Arlong's rows); those two now only keep
# their plain nickname.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKey = @'
b'F\x16v2\xbb\xe1"\x8c\x04s\xb2*\x957K\xe1z\x02>M\x7f\x07\xda\xa8I\x86:\' :\xa8\x0f\xcb\x97\xaf\x8f4\xd6?&\x0f\xef\xf8%\x96{EGE\xe7F\x97e\x91\xdbPmC\xc1w~\x96\xb0A\xcef\xb7\x12%\r4NLS\x06>\x1c\xbc\xb8U\x824\x8a\xe8\xafPau]\'^?h\xa6\x17K\xcb!\xe4\x04CfS\xe7\x02\x16\xd1EsF\r>\x86\tP\xaa-\xa5\x89\r\xb8\xb2G\xe0h\xa62\xe6r|\x8b\xf5\xbd\x93bzK\xc1\xec\x19\xe9\x16/\xbd)Ef\xc9\xad\x84\x8a\xbf\xaeyf\r\x16\xd7\xffK\xb3h\xd9+*\xaf\xb3\xbeB\xd0\x10*\xd7\xd5\x19.r\xbe.z\xf4lSq9\xb5A\xdej\xce=U\xad\xf7_mz3\xe6\xe5R)\x10\x0e\xee\xe2\tP\xa5`\xcd\x068\x08\x98\x89\xb4(\x82e\x83\x1b\xf4\x06\xb8/\x9f?\xf9\xe53VzK_\xae{t@P\xce4\xc9\xcb\nR\xf6\xe8\xd79\xa5j\xde\r\x9cI'
'@

$newToken = @'
b'gAAAAABlNulwl57UPR2nFW4ep293W7eF2DM1VWmeJiQCvkXBOnocaVTNnzCI1E39QUgKxfy8IOBbciCat4rEsFEaRGpsLIr_-Q=='
'@

# Luffy's row (row 2) gets the new symmetric key + coordinates token
$ws.Cells.Item(2, 2).Value = $newKey
$ws.Cells.Item(2, 3).Value = $newToken

# Remove the old key/token pairs from Sanji's row (row 5) and Arlong's row (row 9)
$ws.Cells.Item(5, 2).ClearContents()
$ws.Cells.Item(5, 3).ClearContents()
$ws.Cells.Item(9, 2).ClearContents()
$ws.Cells.Item(9, 3).ClearContents()

# Column B widens slightly to fit the new content
$ws.Columns.Item(2).ColumnWidth = 117.6667
